# au cafe et au tabac
# Adds a new vocabulary lesson block ("quatrième leçon") to Sheet1,
# covering ordering drinks at the café and buying things at the tabac.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the two B-cells that were missing on existing rows ---
$ws.Cells.Item(674, 2).Value = "slices (of bread)"
$ws.Cells.Item(675, 2).Value = "bơ"

# --- New rows 679-691: English / Français / IPA-Français / Gender ---

# Row 679: lesson title (no English / gender)
$ws.Cells.Item(679, 3).Value = "quatrième leçon"
$ws.Cells.Item(679, 4).Value = "ka.tʁi.jɛm"

# Row 680
$ws.Cells.Item(680, 2).Value = "let's order"
$ws.Cells.Item(680, 3).Value = "commandons"
$ws.Cells.Item(680, 4).Value = "kɔ.mɑ̃.do"

# Row 681
$ws.Cells.Item(681, 2).Value = "black"
$ws.Cells.Item(681, 3).Value = "noirs"
$ws.Cells.Item(681, 4).Value = "nwaʁ"

# Row 682
$ws.Cells.Item(682, 2).Value = "beer"
$ws.Cells.Item(682, 3).Value = "bières"
$ws.Cells.Item(682, 4).Value = "bjɛʁ"
$ws.Cells.Item(682, 5).Value = "F"

# Row 683
$ws.Cells.Item(683, 2).Value = "german"
$ws.Cells.Item(683, 3).Value = "allemand / allemandes"
$ws.Cells.Item(683, 4).Value = "al.mɑ̃"

# Row 684
$ws.Cells.Item(684, 2).Value = "glass"
$ws.Cells.Item(684, 3).Value = "verre"
$ws.Cells.Item(684, 4).Value = "vɛʁ"

# Row 685
$ws.Cells.Item(685, 2).Value = "white"
$ws.Cells.Item(685, 3).Value = "blanc"
$ws.Cells.Item(685, 4).Value = "blɑ̃"

# Row 686
$ws.Cells.Item(686, 2).Value = "packet"
$ws.Cells.Item(686, 3).Value = "paquets"
$ws.Cells.Item(686, 4).Value = "pa.kɛ"

# Row 687
$ws.Cells.Item(687, 2).Value = "brown cigarettes"
$ws.Cells.Item(687, 3).Value = "cigarettes brunes"
$ws.Cells.Item(687, 4).Value = "si.ga.ʁɛt bʁyn"
$ws.Cells.Item(687, 5).Value = "F"

# Row 688
$ws.Cells.Item(688, 2).Value = "tobacco shop"
$ws.Cells.Item(688, 3).Value = "tabac"
$ws.Cells.Item(688, 4).Value = "ta.ba"

# Row 689
$ws.Cells.Item(689, 2).Value = "dutch cigar"
$ws.Cells.Item(689, 3).Value = "cigare hollandais"
$ws.Cells.Item(689, 4).Value = "si.ga ɔ.lɑ̃.dɛ"

# Row 690
$ws.Cells.Item(690, 2).Value = "red lighter"
$ws.Cells.Item(690, 3).Value = "briquet rouge"
$ws.Cells.Item(690, 4).Value = "bʁi.kɛt ʁuʒ"

# Row 691
$ws.Cells.Item(691, 2).Value = "im sorry"
$ws.Cells.Item(691, 3).Value = "je suis désolé"
$ws.Cells.Item(691, 4).Value = "ʒə sɥi de.zɔ.le"

# --- Restore the view/selection to match where the author ended up
#     after entering the new rows ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 680
$ws.Range("D697").Select()

